$d = $word.ActiveDocument

function Set-ParagraphXml($range, $innerWordXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerWordXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ------------------------------------------------------------------
# 1) Table row "Do chinh xac" -> value cell: "97%" becomes "9" + "6" + "%"
#    (the "6" run is a fresh run with no w:lang, "9" keeps the original run).
# ------------------------------------------------------------------
$table = $d.Tables.Item(1)
$accuracyCell = $table.Cell(2, 2)
$accuracyPara = $accuracyCell.Range.Paragraphs.Item(1)
$accuracyContent = $d.Range($accuracyPara.Range.Start, $accuracyPara.Range.End - 1)

$accuracyXml = '<w:p w:rsidR="00F95C2E" w:rsidRPr="00D01A82" w:rsidRDefault="009363FA" w:rsidP="00BA71FB">' +
  '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="vi-VN"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="vi-VN"/></w:rPr><w:t>9</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>6</w:t></w:r>' +
  '<w:r w:rsidR="00F95C2E" w:rsidRPr="00D01A82"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="vi-VN"/></w:rPr><w:t>%</w:t></w:r>' +
  '</w:p>'

Set-ParagraphXml $accuracyContent $accuracyXml

# ------------------------------------------------------------------
# 2) Table row "So cau hoi rut sai" -> value cell: "1" becomes "2",
#    and the vi-VN language mark is dropped from both the paragraph
#    mark run properties and the run itself.
# ------------------------------------------------------------------
$wrongCountCell = $table.Cell(5, 2)
$wrongCountPara = $wrongCountCell.Range.Paragraphs.Item(1)
$wrongCountContent = $d.Range($wrongCountPara.Range.Start, $wrongCountPara.Range.End - 1)

$wrongCountXml = '<w:p w:rsidR="00F95C2E" w:rsidRPr="00D01A82" w:rsidRDefault="009363FA" w:rsidP="00BA71FB">' +
  '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>2</w:t></w:r>' +
  '</w:p>'

Set-ParagraphXml $wrongCountContent $wrongCountXml

# ------------------------------------------------------------------
# 3) Final narrative paragraph: append a "." as its own run right
#    after the existing "Nguyen nhan: ..." run.
# ------------------------------------------------------------------
$reasonPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Nguy*nhan*") {
        $reasonPara = $candidate
    }
}
$reasonContent = $d.Range($reasonPara.Range.Start, $reasonPara.Range.End - 1)

$reasonXml = '<w:p w:rsidR="00F917C8" w:rsidRPr="009363FA" w:rsidRDefault="009363FA">' +
  '<w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>Nguyên nhân: Thiếu luật, Luật không bao phủ. Gán nhãn sai ngữ nghĩa</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>.</w:t></w:r>' +
  '</w:p>'

Set-ParagraphXml $reasonContent $reasonXml
